# Apply "added new db sheets, new connector symbols" edit.
$wb = $excel.ActiveWorkbook

$wsAudio   = $wb.Worksheets.Item("audio")
$wsDC      = $wb.Worksheets.Item("DC power")
$wsBattery = $wb.Worksheets.Item("battery")

# --- Common header row used on every sheet ---
$headers = @("TPN","Description","Value","Tolerance","Power Rating","Package","Pulse-Rated","Library Ref","Footprint Ref","Manufacturer 1","Manufacturer 1 PN","Supplier 1","Supplier 1 PN","Manufacturer 2","Manufacturer 2 PN","Supplier 2","Supplier 2 PN")

# --- Sheet "DC power" (sheet2): brand new header row ---
for ($i = 0; $i -lt $headers.Length; $i++) {
    $wsDC.Cells.Item(1, $i + 1).Value = $headers[$i]
}
$wsDC.Range("A1:Q1").Font.Bold = $true
$wsDC.Range("C1").NumberFormat = "@"
$wsDC.Range("F1").NumberFormat = "@"

# --- Sheet "battery" (sheet3): brand new header row ---
for ($i = 0; $i -lt $headers.Length; $i++) {
    $wsBattery.Cells.Item(1, $i + 1).Value = $headers[$i]
}
$wsBattery.Range("A1:Q1").Font.Bold = $true
$wsBattery.Range("C1").NumberFormat = "@"
$wsBattery.Range("F1").NumberFormat = "@"

# --- New connector symbols (order drives the shared-string table layout) ---
$wsAudio.Range("H2").Value = "STEREO_JACK"
$wsAudio.Range("H3").Value = "MONO_JACK"
$wsAudio.Range("I2").Value = "STEREO_JACK_HDR"
$wsAudio.Range("I3").Value = "MONO_JACK_HDR"

$wsBattery.Range("I2").Value = "PWR_JACK_HDR"
$wsDC.Range("I2").Value = "PWR_JACK_SW_HDR"

$wsBattery.Range("H2").Value = "9V_BATT"
$wsDC.Range("H2").Value = "PWR_JACK_SW"

# New column widths on the audio sheet for the header/value columns
# (engine quantizes ColumnWidth to whole pixels, so feed it the value
# that lands on the closest achievable width)
$wsAudio.Columns.Item(8).ColumnWidth = 32.83
$wsAudio.Columns.Item(9).ColumnWidth = 17

# --- Final cursor position on each sheet ---
[void]$wsDC.Range("H3").Select()
[void]$wsBattery.Range("H3").Select()
[void]$wsAudio.Range("I19").Select()
